$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values per the diff
$ws.Range("U2").Value = 50
$ws.Range("V2").Value = 40
$ws.Range("X2").Value = 30

# Update the selection shown in the sheet view
$ws.Range("O1:O1048576").Select() | Out-Null
